# "Generate Report for Archive"
#
# The localization-status report was regenerated: the outstanding status
# text moved from "Ready for handoff" to "In Translation", and the
# "Status" column(s) that hold that text were re-sized (narrower, since
# "In Translation" renders slightly differently than "Ready for handoff").

$wb = $excel.ActiveWorkbook

# --- 1. Update the status text everywhere it appears -----------------
# Overview sheet: columns E ("zh-cn") and F ("de-de") hold the current
# status for each locale.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

# Per-locale sheets: column C ("Status") holds the same status value.
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# --- 2. Re-size the columns that hold the status text -----------------
# New width is narrower than before now that the status text changed.
$newWidth = 12.5

$overview.Columns.Item(5).ColumnWidth = $newWidth
$overview.Columns.Item(6).ColumnWidth = $newWidth

$zhcn.Columns.Item(3).ColumnWidth = $newWidth
$dede.Columns.Item(3).ColumnWidth = $newWidth
